$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string (row 1, column A)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 21:22"

# Update Murcia's stats (row 30: B=Casos totales, C=Casos activos, D=Recuperados, E=Muertes)
$ws.Range("B30").Value = 1283
$ws.Range("C30").Value = 193
$ws.Range("D30").Value = 1012
$ws.Range("E30").Value = 78
